$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") for data rows 2-372 all change from 45178 to 45179.
$ws.Range("C2:C372").Value = 45179
